# chore: update Sheets via scheduled runner
# Refresh market-board price/profit columns (H:N) on each crafting-leve
# sheet with the latest snapshot values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H15").Value = 560.4737
$ws.Range("I15").Value = 560.4737
$ws.Range("K15").Value = 1681.4211
$ws.Range("M15").Value = -1512.4211

$ws.Range("H38").Value = 3583.2
$ws.Range("I38").Value = 38.666668
$ws.Range("K38").Value = 116.000004
$ws.Range("M38").Value = 255.999996

$ws.Range("H41").Value = 361.44446
$ws.Range("I41").Value = 414
$ws.Range("K41").Value = 414
$ws.Range("M41").Value = 26

$ws.Range("H43").Value = 12499
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10138

$ws.Range("H47").Value = 4075
$ws.Range("J47").Value = 4075
$ws.Range("L47").Value = 4075
$ws.Range("N47").Value = -6019

$ws.Range("H49").Value = 1000
$ws.Range("I49").Value = 1000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 3000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -2864
$ws.Range("N49").ClearContents()

$ws.Range("H58").Value = 3270.1667
$ws.Range("I58").Value = 68
$ws.Range("J58").Value = 4871.25
$ws.Range("K58").Value = 204
$ws.Range("L58").Value = 14613.75
$ws.Range("M58").Value = -54
$ws.Range("N58").Value = -14913.75

$ws.Range("H88").Value = 5527.857
$ws.Range("I88").Value = 3347.5
$ws.Range("K88").Value = 3347.5
$ws.Range("M88").Value = -2941.5

$ws.Range("H91").Value = 5527.857
$ws.Range("I91").Value = 3347.5
$ws.Range("K91").Value = 3347.5
$ws.Range("M91").Value = -1943.5

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H112").Value = 2600
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 3000
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 9000
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -11216

$ws.Range("H132").Value = 6013.3335
$ws.Range("I132").Value = 6013.3335
$ws.Range("K132").Value = 18040.0005
$ws.Range("M132").Value = -15510.0005

$ws.Range("H138").Value = 3395.2727
$ws.Range("I138").Value = 1724.6666
$ws.Range("K138").Value = 5173.9998
$ws.Range("M138").Value = -33.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 3000
$ws.Range("J13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3288

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3694.8572
$ws.Range("I20").Value = 3694.8572
$ws.Range("K20").Value = 3694.8572
$ws.Range("M20").Value = -3447.8572

$ws.Range("H76").Value = 25396.8
$ws.Range("I76").Value = 26999
$ws.Range("J76").Value = 24328.666
$ws.Range("K76").Value = 26999
$ws.Range("L76").Value = 24328.666
$ws.Range("M76").Value = -26684
$ws.Range("N76").Value = -24958.666

$ws.Range("H79").Value = 25396.8
$ws.Range("I79").Value = 26999
$ws.Range("J79").Value = 24328.666
$ws.Range("K79").Value = 26999
$ws.Range("L79").Value = 24328.666
$ws.Range("M79").Value = -25907
$ws.Range("N79").Value = -26512.666

$ws.Range("H86").Value = 1031.7778
$ws.Range("I86").Value = 1214.1428
$ws.Range("J86").Value = 393.5
$ws.Range("K86").Value = 1214.1428
$ws.Range("L86").Value = 393.5
$ws.Range("M86").Value = -91.14280000000008
$ws.Range("N86").Value = -2639.5

$ws.Range("H89").Value = 1031.7778
$ws.Range("I89").Value = 1214.1428
$ws.Range("J89").Value = 393.5
$ws.Range("K89").Value = 6070.714
$ws.Range("L89").Value = 1967.5
$ws.Range("M89").Value = -454.7139999999999
$ws.Range("N89").Value = -13199.5

$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 556
$ws.Range("J34").Value = 1132.3334
$ws.Range("L34").Value = 3397.0002
$ws.Range("N34").Value = -3565.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 19
$ws.Range("I2").Value = 17.714285
$ws.Range("J2").Value = 20.5
$ws.Range("K2").Value = 17.714285
$ws.Range("L2").Value = 20.5
$ws.Range("M2").Value = 95.285715
$ws.Range("N2").Value = -246.5

$ws.Range("H6").Value = 1468.75
$ws.Range("I6").Value = 1468.75
$ws.Range("K6").Value = 1468.75
$ws.Range("M6").Value = -1355.75

$ws.Range("H16").Value = 1468.75
$ws.Range("I16").Value = 1468.75
$ws.Range("K16").Value = 1468.75
$ws.Range("M16").Value = -1218.75

$ws.Range("H17").Value = 783.3333
$ws.Range("J17").Value = 150
$ws.Range("L17").Value = 150
$ws.Range("N17").Value = -486

$ws.Range("H31").Value = 952.8570999999999
$ws.Range("I31").Value = 952.8570999999999
$ws.Range("K31").Value = 952.8570999999999
$ws.Range("M31").Value = -660.8570999999999

$ws.Range("H37").Value = 952.8570999999999
$ws.Range("I37").Value = 952.8570999999999
$ws.Range("K37").Value = 952.8570999999999
$ws.Range("M37").Value = -675.8570999999999

$ws.Range("H80").Value = 4385
$ws.Range("I80").Value = 3262
$ws.Range("K80").Value = 3262
$ws.Range("M80").Value = -2264

$ws.Range("H83").Value = 4385
$ws.Range("I83").Value = 3262
$ws.Range("K83").Value = 16310
$ws.Range("M83").Value = -11318

$ws.Range("H122").Value = 1638.2
$ws.Range("I122").Value = 1638.2
$ws.Range("K122").Value = 4914.6
$ws.Range("M122").Value = -2464.6

$ws.Range("H132").Value = 7304.2354
$ws.Range("I132").Value = 5507.364
$ws.Range("J132").Value = 10598.5
$ws.Range("K132").Value = 16522.092
$ws.Range("L132").Value = 31795.5
$ws.Range("M132").Value = -13992.092
$ws.Range("N132").Value = -36855.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10124.5
$ws.Range("I7").Value = 10124.5
$ws.Range("K7").Value = 10124.5
$ws.Range("M7").Value = -10012.5

$ws.Range("H46").Value = 5650
$ws.Range("J46").Value = 5630.4346
$ws.Range("L46").Value = 5630.4346
$ws.Range("N46").Value = -6006.4346

$ws.Range("H53").Value = 7000
$ws.Range("I53").Value = 7000
$ws.Range("K53").Value = 7000
$ws.Range("M53").Value = -6482

$ws.Range("H68").Value = 4899.25
$ws.Range("I68").Value = 4899.25
$ws.Range("K68").Value = 4899.25
$ws.Range("M68").Value = -4150.25

$ws.Range("H71").Value = 4899.25
$ws.Range("I71").Value = 4899.25
$ws.Range("K71").Value = 24496.25
$ws.Range("M71").Value = -20752.25

$ws.Range("H100").Value = 4541.5557
$ws.Range("I100").Value = 3410.8572
$ws.Range("K100").Value = 3410.8572
$ws.Range("M100").Value = -2869.8572

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H126").Value = 10124.5
$ws.Range("I126").Value = 10124.5
$ws.Range("K126").Value = 30373.5
$ws.Range("M126").Value = -27903.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 22299.6
$ws.Range("I62").Value = 2749.5
$ws.Range("J62").Value = 35333
$ws.Range("K62").Value = 2749.5
$ws.Range("L62").Value = 35333
$ws.Range("M62").Value = -2125.5
$ws.Range("N62").Value = -36581

$ws.Range("H65").Value = 22299.6
$ws.Range("I65").Value = 2749.5
$ws.Range("J65").Value = 35333
$ws.Range("K65").Value = 13747.5
$ws.Range("L65").Value = 176665
$ws.Range("M65").Value = -10627.5
$ws.Range("N65").Value = -182905

$ws.Range("H81").Value = 6770
$ws.Range("I81").Value = 1925
$ws.Range("K81").Value = 3850
$ws.Range("M81").Value = -2789

$ws.Range("H84").Value = 6770
$ws.Range("I84").Value = 1925
$ws.Range("K84").Value = 19250
$ws.Range("M84").Value = -13946

$ws.Range("H96").Value = 1243.5
$ws.Range("I96").Value = 1243.5
$ws.Range("K96").Value = 1243.5
$ws.Range("M96").Value = 129.5

$ws.Range("H122").Value = 5334.6665
$ws.Range("I122").Value = 5334.6665
$ws.Range("K122").Value = 16003.9995
$ws.Range("M122").Value = -13553.9995
